$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1): update F2, F3, F4
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 1013
$ws1.Range("F3").Value = 2090
$ws1.Range("F4").Value = 458

# Sheet "全部类型" (sheet4): update F4, F5, F6
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 1013
$ws4.Range("F5").Value = 2090
$ws4.Range("F6").Value = 458
